# Generate Report for Handoff
#
# A new handoff was generated for 01815df5-c612-4267-8e3e-93304d033164.md
# (row 5 of each sheet). Refresh the "Latest Handoff Date(time)" stamps for
# that file on the Overview sheet and on each per-locale (zh-cn / de-de)
# detail sheet.

$wb = $excel.ActiveWorkbook

# Overview sheet: column D is "Latest Handoff Date", row 5 is 01815df5...md
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D5").Value = "2016-40-13 12:40:13"

# zh-cn sheet: column E is "Latest Handoff Datetime", row 5 is 01815df5...md
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E5").Value = "2016-03-13 12:40:09"

# de-de sheet: column E is "Latest Handoff Datetime", row 5 is 01815df5...md
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E5").Value = "2016-03-13 12:40:13"
